# Applies the docx edit described by the diff:
#  - row 1 ("99.99") -> "0M"
#  - row 2 ("0")      -> "0M"
#  - row 3 ("37")     -> "0M", followed by 10 new single-cell rows
#  - the big tab-separated row ("100\t0.00003\t...\t100.0") -> "99.99"
#  - the big tab-separated row ("3\t0.00004\t...\t100.0")   -> "0"
#  - the trailing empty row -> "37"

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Rows 1-3: simple text replacements ---
$t.Cell(1,1).Range.Text = "0M"
$t.Cell(2,1).Range.Text = "0M"
$t.Cell(3,1).Range.Text = "0M"

# --- Insert 10 new rows right after row 3, each holding one value ---
$newValues = @("103","0.00003","0.00006","0.00005","0.00001","0.00004","0.00004","0.00005","0.00415","100.0")

$afterIndex = 3
foreach ($val in $newValues) {
    $refRow = $t.Rows.Item($afterIndex + 1)
    $newRow = $t.Rows.Add($refRow)
    $newRow.Cells.Item(1).Range.Text = $val
    $afterIndex = $afterIndex + 1
}

# --- The two multi-run (tab separated) rows collapse to a single value ---
# These were originally rows 34 and 35 (out of 36); after inserting 10 new
# rows earlier (right after row 3) everything from row 4 onward shifted down
# by 10, so they now sit at rows 44 and 45.
$t.Cell(44,1).Range.Text = "99.99"
$t.Cell(45,1).Range.Text = "0"

# --- Trailing empty row (was row 36, now row 46) gains the text "37" ---
$t.Cell(46,1).Range.Text = "37"

Write-Output ("Final row count=" + $t.Rows.Count)
